$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: "04/01/2018" -> "09/01/2018"
#    Appears on the Slide Master, every Slide Layout and the Notes
#    Master (ppPlaceholderDate = 16).
# ---------------------------------------------------------------------
$ppPlaceholderDate = 16

function Update-DateePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Placeholders.Count; $i++) {
        $shp = $shapes.Placeholders.Item($i)
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = "09/01/2018"
        }
    }
}

$master = $p.SlideMaster
Update-DateePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateePlaceholder $layout.Shapes
}

$notesMaster = $p.NotesMaster
Update-DateePlaceholder $notesMaster.Shapes

# ---------------------------------------------------------------------
# 2) Slide 1: reposition + relabel the rotated axis-title textboxes.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$ebv1 = $s1.Shapes.Item(2)   # "TextBox 4" -> "EBVhost"
$ebv1.Left = 296.8495330840551
$ebv1.Top = 216.3436203022441
$ebv1.TextFrame.TextRange.Text = "EBVhost (mag)"

$tmax1 = $s1.Shapes.Item(3)  # "TextBox 6" -> "Tmax"
$tmax1.Left = 296.8448028596063
$tmax1.Top = 124.1355934161811
$tmax1.TextFrame.TextRange.Text = "Tmax (days)"

# "TextBox 7" holds "DM " + "(mag)" as two separate runs; merge them into
# a single run "DM (mag)" (force a real change first so the merge happens
# even though the visible text ends up the same).
$dm1 = $s1.Shapes.Item(4)
$dm1.TextFrame.TextRange.Text = "DM (mag)__tmp__"
$dm1.TextFrame.TextRange.Text = "DM (mag)"

# ---------------------------------------------------------------------
# 3) Slide 2: same run-merge fix for its "DM (mag)" textbox.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$dm2 = $s2.Shapes.Item(5)    # "TextBox 6" -> "DM " + "(mag)"
$dm2.TextFrame.TextRange.Text = "DM (mag)__tmp__"
$dm2.TextFrame.TextRange.Text = "DM (mag)"
